# Architecture_Diagram.pptx edit: rename two boxes on slide 1
#   "rest api"                      -> "Producer"
#   "Spark Streaming\nApplication"  -> "Consumer"
#
# The shapes are rewritten by grabbing the *whole* TextRange (via
# Characters so the complete run/line-break span is addressed), setting
# its .Text (which collapses every run + <a:br> into a single run), and
# then re-asserting the run-level formatting (size/underline/strike/
# shadow/typeface) that PowerPoint normalizes onto freshly retyped text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Set-BoxText($shape, $newText) {
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Characters(1, $tr.Length)
    $full.Text = $newText
    $full.Font.Size = 18
    $full.Font.Underline = $false
    $full.Font.Strikethrough = $false
    $full.Font.Shadow = $false
    $full.Font.Name = "+mn-lt"
}

# Shape 1 ("직사각형 4"): "rest api" -> "Producer"
Set-BoxText $s.Shapes.Item(1) "Producer"

# Shape 3 ("직사각형 6"): "Spark Streaming" / "Application" -> "Consumer"
Set-BoxText $s.Shapes.Item(3) "Consumer"
